# daily auto push: 2026-01-26 22:37 UTC
# A new reading (2026/01/27, time=5, ranking=23) was appended into the
# time-series right after the existing "2026/01/27" row, pushing every
# following row down by one (old row 705 -> new row 706, ..., old row
# 746 -> new row 747).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at 705; everything from the old row 705
# downward shifts to row+1 automatically.
$ws.Rows(705).Insert()

# Column A holds dates formatted/stored as literal text (e.g. "2026/01/27"),
# not real Excel date serials. Prefixing with an apostrophe forces the
# value to be stored as text instead of being auto-parsed into a date
# serial number; resetting the style back to "Normal" afterwards clears
# the quote-prefix formatting flag so the new cell matches the plain,
# unstyled look of every other data row.
$ws.Range("A705").Value = "'2026/01/27"
$ws.Range("A705").Style = "Normal"

$ws.Range("B705").Value = "火"
$ws.Range("C705").Value = 5
$ws.Range("D705").Value = 23
